$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect()

# Update the confidential disclaimer date (2021-05-10 -> 2021-05-11)
$ws.Range("A80").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-11 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-77
$ws.Cells.Item(2,4).Value = 0.06266477728102028
$ws.Cells.Item(2,5).Value = -0.00741032715806067
$ws.Cells.Item(3,4).Value = 0.03752676488303065
$ws.Cells.Item(3,5).Value = 0.01047488003410146
$ws.Cells.Item(4,4).Value = 0.03198083152796676
$ws.Cells.Item(4,5).Value = -0.003843353022089269
$ws.Cells.Item(5,4).Value = 0.02923662311018145
$ws.Cells.Item(5,5).Value = -0.002168429437536035
$ws.Cells.Item(6,4).Value = 0.02695572260708715
$ws.Cells.Item(6,5).Value = -0.009464383113341301
$ws.Cells.Item(7,4).Value = 0.02528374443817805
$ws.Cells.Item(7,5).Value = -0.01662324773601298
$ws.Cells.Item(8,4).Value = 0.1894594936245624
$ws.Cells.Item(8,5).Value = -0.0242741551642075
$ws.Cells.Item(9,4).Value = 0.02470030658774925
$ws.Cells.Item(9,5).Value = -0.008163505021436701
$ws.Cells.Item(10,4).Value = 0.02266997737032732
$ws.Cells.Item(10,5).Value = -0.004939347715551601
$ws.Cells.Item(11,4).Value = 0.02213103934922697
$ws.Cells.Item(11,5).Value = -0.007589447054571763
$ws.Cells.Item(12,4).Value = 0.02023232768868224
$ws.Cells.Item(12,5).Value = -0.01427021161150321
$ws.Cells.Item(13,4).Value = 0.02011353078747941
$ws.Cells.Item(13,5).Value = -0.01141498216409031
$ws.Cells.Item(14,4).Value = 0.01729918122901691
$ws.Cells.Item(14,5).Value = -0.006207674943566666
$ws.Cells.Item(15,4).Value = 0.01593148779615901
$ws.Cells.Item(15,5).Value = 0.01086763924162781
$ws.Cells.Item(16,4).Value = 0.01460378539934959
$ws.Cells.Item(16,5).Value = 0.001932989690721643
$ws.Cells.Item(17,4).Value = 0.01435490001029495
$ws.Cells.Item(17,5).Value = -0.009018605311745365
$ws.Cells.Item(18,4).Value = 0.01438109021821689
$ws.Cells.Item(18,5).Value = -0.007292802617230087
$ws.Cells.Item(19,4).Value = 0.01319574806836043
$ws.Cells.Item(19,5).Value = 0.001830244795241409
$ws.Cells.Item(20,4).Value = 0.01373997902070249
$ws.Cells.Item(20,5).Value = -0.03179929689996808
$ws.Cells.Item(21,4).Value = 0.01266527873724681
$ws.Cells.Item(21,5).Value = -0.01133925835121075
$ws.Cells.Item(22,4).Value = 0.01321421451735928
$ws.Cells.Item(22,5).Value = -0.01144381345723633
$ws.Cells.Item(23,4).Value = 0.01149044403607759
$ws.Cells.Item(23,5).Value = 0.003466724900365925
$ws.Cells.Item(24,4).Value = 0.013201433068583
$ws.Cells.Item(24,5).Value = -0.02183468364655849
$ws.Cells.Item(25,4).Value = 0.01203683136781443
$ws.Cells.Item(25,5).Value = -0.03066369606003749
$ws.Cells.Item(26,4).Value = 0.008875657588876608
$ws.Cells.Item(26,5).Value = -0.002102659245516336
$ws.Cells.Item(27,4).Value = 0.009197154153187892
$ws.Cells.Item(27,5).Value = 0.01019694773637969
$ws.Cells.Item(28,4).Value = 0.01028757640277635
$ws.Cells.Item(28,5).Value = -0.002210433244916099
$ws.Cells.Item(29,4).Value = 0.009799803383679201
$ws.Cells.Item(29,5).Value = 0.005509079779636883
$ws.Cells.Item(30,4).Value = 0.009820543832767084
$ws.Cells.Item(30,5).Value = -0.01189715745768116
$ws.Cells.Item(31,4).Value = 0.008558199335068259
$ws.Cells.Item(31,5).Value = -0.01661604430945152
$ws.Cells.Item(32,4).Value = 0.0103720280612552
$ws.Cells.Item(32,5).Value = -0.0001738828029907991
$ws.Cells.Item(33,4).Value = 0.009567737753903774
$ws.Cells.Item(33,5).Value = -0.01646505376344098
$ws.Cells.Item(34,4).Value = 0.009041973250560565
$ws.Cells.Item(34,5).Value = -0.01074485521762869
$ws.Cells.Item(35,4).Value = 0.009372801056616493
$ws.Cells.Item(35,5).Value = -0.01125240525391114
$ws.Cells.Item(36,4).Value = 0.008393883225186725
$ws.Cells.Item(36,5).Value = -0.004624180258954169
$ws.Cells.Item(37,4).Value = 0.008669821373921213
$ws.Cells.Item(37,5).Value = -0.01126034459367797
$ws.Cells.Item(38,4).Value = 0.007398812151745217
$ws.Cells.Item(38,5).Value = -0.01882233244308773
$ws.Cells.Item(39,4).Value = 0.008751606963330156
$ws.Cells.Item(39,5).Value = -0.01279478173607629
$ws.Cells.Item(40,4).Value = 0.008162209664392166
$ws.Cells.Item(40,5).Value = -0.02619330108606366
$ws.Cells.Item(41,4).Value = 0.006839251302251239
$ws.Cells.Item(41,5).Value = 0.009103416647557871
$ws.Cells.Item(42,4).Value = 0.007077472415026282
$ws.Cells.Item(42,5).Value = -0.004420660772452178
$ws.Cells.Item(43,4).Value = 0.008170325492296123
$ws.Cells.Item(43,5).Value = -0.03069245165315038
$ws.Cells.Item(44,4).Value = 0.007486949258644214
$ws.Cells.Item(44,5).Value = 0.005278592375366431
$ws.Cells.Item(45,4).Value = 0.007268488422504889
$ws.Cells.Item(45,5).Value = 0.001143547586683002
$ws.Cells.Item(46,4).Value = 0.007949277000883073
$ws.Cells.Item(46,5).Value = -0.01444128787878773
$ws.Cells.Item(47,4).Value = 0.007355527736257455
$ws.Cells.Item(47,5).Value = -0.009082768325444524
$ws.Cells.Item(48,4).Value = 0.00719964110946457
$ws.Cells.Item(48,5).Value = -0.01678356713426865
$ws.Cells.Item(49,4).Value = 0.006722022676971881
$ws.Cells.Item(49,5).Value = -0.006561679790026198
$ws.Cells.Item(50,4).Value = 0.007405006841642922
$ws.Cells.Item(50,5).Value = -0.01723937099592299
$ws.Cells.Item(51,4).Value = 0.006658389881377099
$ws.Cells.Item(51,5).Value = -0.01364329582457435
$ws.Cells.Item(52,4).Value = 0.006723590952895351
$ws.Cells.Item(52,5).Value = -0.01556942095749025
$ws.Cells.Item(53,4).Value = 0.00534350814024204
$ws.Cells.Item(53,5).Value = -0.02003081664098605
$ws.Cells.Item(54,4).Value = 0.006108748377099065
$ws.Cells.Item(54,5).Value = -0.00005134524543026764
$ws.Cells.Item(55,4).Value = 0.006142740757740271
$ws.Cells.Item(55,5).Value = -0.04643370033508853
$ws.Cells.Item(56,4).Value = 0.005733020831354202
$ws.Cells.Item(56,5).Value = -0.01060593004743393
$ws.Cells.Item(57,4).Value = 0.006779029506790004
$ws.Cells.Item(57,5).Value = -0.008328320917965981
$ws.Cells.Item(58,4).Value = 0.005547932906866311
$ws.Cells.Item(58,5).Value = -0.004748982360922693
$ws.Cells.Item(59,4).Value = 0.00530210565586244
$ws.Cells.Item(59,5).Value = -0.003460668175162884
$ws.Cells.Item(60,4).Value = 0.005009073299562131
$ws.Cells.Item(60,5).Value = -0.006386975579210863
$ws.Cells.Item(61,4).Value = 0.004892393570855988
$ws.Cells.Item(61,5).Value = -0.008021861777150918
$ws.Cells.Item(62,4).Value = 0.004896157433072316
$ws.Cells.Item(62,5).Value = 0.02818705957719425
$ws.Cells.Item(63,4).Value = 0.004293939478459962
$ws.Cells.Item(63,5).Value = -0.01964937910883857
$ws.Cells.Item(64,4).Value = 0.004058698089939511
$ws.Cells.Item(64,5).Value = -0.006800618238021583
$ws.Cells.Item(65,4).Value = 0.003867838910053252
$ws.Cells.Item(65,5).Value = -0.02201678627904147
$ws.Cells.Item(66,4).Value = 0.003708345248636386
$ws.Cells.Item(66,5).Value = 0.00786602385181423
$ws.Cells.Item(67,4).Value = 0.003849960564525698
$ws.Cells.Item(67,5).Value = -0.01643651472565089
$ws.Cells.Item(68,4).Value = 0.00363765521138599
$ws.Cells.Item(68,5).Value = -0.01173731690755653
$ws.Cells.Item(69,4).Value = 0.003685213178765208
$ws.Cells.Item(69,5).Value = -0.03023597250888366
$ws.Cells.Item(70,4).Value = 0.002951652115562269
$ws.Cells.Item(70,5).Value = -0.007013442431326622
$ws.Cells.Item(71,4).Value = 0.002865592974261871
$ws.Cells.Item(71,5).Value = 0.02039978656158925
$ws.Cells.Item(72,4).Value = 0.002237263225523747
$ws.Cells.Item(72,5).Value = 0.00283896745702128
$ws.Cells.Item(73,4).Value = 0.001908160522983637
$ws.Cells.Item(73,5).Value = 0.01723889950481827
$ws.Cells.Item(74,4).Value = 0.001879500280482228
$ws.Cells.Item(74,5).Value = 0.01211982143602164
$ws.Cells.Item(75,4).Value = 0.00144195129783419
$ws.Cells.Item(75,5).Value = -0.01484583174724019
$ws.Cells.Item(76,4).Value = 0.001664254409986016
$ws.Cells.Item(76,5).Value = -0.00381643422540523
$ws.Cells.Item(77,4).Value = 0.9999999999999999
$ws.Cells.Item(77,5).Value = -0.01075624978147061

$ws.Protect()
